$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 ("time_taken"), matching the style of the other
# header cells (E1) - bold font + border + centered alignment.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"
$excel.CutCopyMode = $false

# Timestamp values for F2:F16 (plain, unstyled data cells like the rest
# of the data rows).
$timestamps = @(
    "2021-10-05 13:38:45.889110",
    "2021-10-05 13:38:45.889117",
    "2021-10-05 13:38:45.889120",
    "2021-10-05 13:38:45.889122",
    "2021-10-05 13:38:45.889124",
    "2021-10-05 13:38:45.889126",
    "2021-10-05 13:38:45.889128",
    "2021-10-05 13:38:45.889130",
    "2021-10-05 13:38:45.889132",
    "2021-10-05 13:38:45.889134",
    "2021-10-05 13:38:45.889135",
    "2021-10-05 13:38:45.889137",
    "2021-10-05 13:38:45.889139",
    "2021-10-05 13:38:45.889141",
    "2021-10-05 13:38:45.889143"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
